$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.004.91'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '2.294.14'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''299.84'
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('D6').Value = '''97.71'
$ws.Range('E6').Value = '  -2.73%  '
$ws.Range('E7').Value = '  +2.28%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '''0.513'
$ws.Range('E9').Value = '  -0.26%  '
$ws.Range('D10').Value = '''36.04'
$ws.Range('E10').Value = '  -1.72%  '
$ws.Range('D11').Value = '''0.0788'
$ws.Range('E11').Value = '  -0.59%  '
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('D13').Value = '''17.62'
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('D14').Value = '''6.83'
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('D15').Value = '2.654.48'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').Value = '2.320.46'
$ws.Range('E16').Value = '  +0.71%  '
$ws.Range('D17').Value = '''0.785'
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('D18').Value = '42.924.72'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').Value = '''12.76'
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').Value = '0.0₃0914'
$ws.Range('E20').Value = '  +0.86%  '
$ws.Range('D21').Value = '''6.12'
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('D22').Value = '''69.03'
$ws.Range('E22').Value = '  +1.59%  '
$ws.Range('D23').Value = '''237.22'
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('D24').Value = '''2.12'
$ws.Range('E24').Value = '  -3.33%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').Value = '''2.43'
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('D27').Value = '''24.93'
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('D28').Value = '''165.61'
$ws.Range('E28').Value = '  -1.46%  '
$ws.Range('D29').Value = '''2.03'
$ws.Range('E29').Value = '  -0.82%  '
$ws.Range('D30').Value = '''9.08'
$ws.Range('E30').Value = '  -0.83%  '
$ws.Range('D31').Value = '''33.06'
$ws.Range('E31').Value = '  -4.34%  '
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('D33').Value = '''5.07'
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('D34').Value = '''4.72'
$ws.Range('E34').Value = '  +2.44%  '
$ws.Range('D35').Value = '''17.84'
$ws.Range('E35').Value = '  +1.17%  '
$ws.Range('D36').Value = '''2.40'
$ws.Range('E36').Value = '  -0.53%  '
$ws.Range('D37').Value = '''0.0694'
$ws.Range('E37').Value = '  +0.26%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').Value = '''1.77'
$ws.Range('E39').Value = '  -1.36%  '
$ws.Range('E40').Value = '  +0.41%  '
$ws.Range('D41').Value = '''2.76'
$ws.Range('E41').Value = '  -2.02%  '
$ws.Range('D42').Value = '2.006.26'
$ws.Range('E42').Value = '  +1.14%  '
$ws.Range('E43').Value = '  -2.45%  '
$ws.Range('E44').Value = '  -1.65%  '
$ws.Range('D45').Value = '''10.25'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('D46').Value = '''17.38'
$ws.Range('E46').Value = '  -2.09%  '
$ws.Range('D47').Value = '''2.82'
$ws.Range('E47').Value = '  -3.02%  '
$ws.Range('D48').Value = '''53.99'
$ws.Range('E48').Value = '  -2.54%  '
$ws.Range('D49').Value = '2.521.92'
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').Value = '''73.31'
$ws.Range('E50').Value = '  +3.48%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = '''1.53'
$ws.Range('E51').Value = '  -1.89%  '
